# Auto-generated edit script: updates odds values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 3
$ws.Range("G3").Value = 2.15
$ws.Range("I3").Value = 3.9
$ws.Range("J3").Value = 3
$ws.Range("L3").Value = 4.75
$ws.Range("W3").Value = 2.2
$ws.Range("X3").Value = 1.62
$ws.Range("Z3").Value = 8.5
$ws.Range("AA3").Value = 10
$ws.Range("AB3").Value = 19
$ws.Range("AJ3").Value = 8
$ws.Range("AK3").Value = 17
$ws.Range("AM3").Value = 41

# Row 7
$ws.Range("J7").Value = 2.35
$ws.Range("L7").Value = 5.2
$ws.Range("Q7").Value = 2.15
$ws.Range("R7").Value = 1.55
$ws.Range("X7").Value = 1.6
$ws.Range("Y7").Value = 5.4
$ws.Range("AC7").Value = 16
$ws.Range("AD7").Value = 37
$ws.Range("AE7").Value = 7.7
$ws.Range("AJ7").Value = 11
$ws.Range("AL7").Value = 16.5
$ws.Range("AM7").Value = 90

# Row 10
$ws.Range("G10").Value = 1.85
$ws.Range("I10").Value = 5
$ws.Range("M10").Value = 1.11
$ws.Range("N10").Value = 6.5
$ws.Range("R10").Value = 1.48
$ws.Range("X10").Value = 1.57
$ws.Range("AK10").Value = 23

# Row 11
$ws.Range("R11").Value = 1.44
$ws.Range("X11").Value = 1.62
$ws.Range("AP11").Value = 2.03
$ws.Range("AQ11").Value = 1.83
$ws.Range("AR11").Value = 4.4

# Row 12
$ws.Range("K12").Value = 2.38
$ws.Range("Q12").Value = 1.57
$ws.Range("R12").Value = 2.35
$ws.Range("S12").Value = 2.38
$ws.Range("T12").Value = 1.53
$ws.Range("U12").Value = 1.3
$ws.Range("V12").Value = 3.4
$ws.Range("W12").Value = 1.5
$ws.Range("AA12").Value = 9.5
$ws.Range("AJ12").Value = 13
$ws.Range("AK12").Value = 17
$ws.Range("AM12").Value = 29
$ws.Range("AR12").Value = 1.9
$ws.Range("AS12").Value = 1.9

# Row 13
$ws.Range("J13").Value = 2.72
$ws.Range("K13").Value = 2.1
$ws.Range("L13").Value = 3.65
$ws.Range("O13").Value = 1.32
$ws.Range("P13").Value = 2.85
$ws.Range("Q13").Value = 1.93
$ws.Range("W13").Value = 1.72
$ws.Range("X13").Value = 1.9
$ws.Range("Y13").Value = 7.7
$ws.Range("Z13").Value = 11
$ws.Range("AC13").Value = 18
$ws.Range("AD13").Value = 27
$ws.Range("AE13").Value = 8.75
$ws.Range("AF13").Value = 6.1
$ws.Range("AH13").Value = 65
$ws.Range("AJ13").Value = 9
$ws.Range("AK13").Value = 16
$ws.Range("AL13").Value = 11.25
$ws.Range("AM13").Value = 40
$ws.Range("AN13").Value = 30
$ws.Range("AO13").Value = 37

# Row 14
$ws.Range("G14").Value = 2.57
$ws.Range("I14").Value = 2.65
$ws.Range("J14").Value = 3.05
$ws.Range("K14").Value = 2.07
$ws.Range("L14").Value = 3.2
$ws.Range("P14").Value = 2.82
$ws.Range("X14").Value = 1.88
$ws.Range("Y14").Value = 8.25
$ws.Range("Z14").Value = 13
$ws.Range("AA14").Value = 9.5
$ws.Range("AB14").Value = 29
$ws.Range("AC14").Value = 21
$ws.Range("AD14").Value = 30
$ws.Range("AF14").Value = 6
$ws.Range("AG14").Value = 14
$ws.Range("AJ14").Value = 8
$ws.Range("AK14").Value = 13
$ws.Range("AL14").Value = 10
$ws.Range("AM14").Value = 30
$ws.Range("AN14").Value = 23
$ws.Range("AO14").Value = 32

# Row 15
$ws.Range("J15").Value = 1.91
$ws.Range("M15").Value = 1.04
$ws.Range("N15").Value = 13
$ws.Range("Q15").Value = 1.67
$ws.Range("U15").Value = 1.33
$ws.Range("V15").Value = 3.25
$ws.Range("W15").Value = 1.91
$ws.Range("X15").Value = 1.8
$ws.Range("AC15").Value = 12
$ws.Range("AJ15").Value = 19

# Row 16
$ws.Range("Q16").Value = 1.8
$ws.Range("R16").Value = 2

# Row 17
$ws.Range("O17").Value = 1.13
$ws.Range("P17").Value = 6
$ws.Range("Q17").Value = 1.44
$ws.Range("R17").Value = 2.7
$ws.Range("AE17").Value = 21
$ws.Range("AI17").Value = 101
$ws.Range("AK17").Value = 11

# Row 18
$ws.Range("G18").Value = 4.75
$ws.Range("H18").Value = 4.75
$ws.Range("I18").Value = 1.57
$ws.Range("M18").Value = 1.01
$ws.Range("N18").Value = 26
$ws.Range("O18").Value = 1.1
$ws.Range("P18").Value = 7
$ws.Range("Q18").Value = 1.36
$ws.Range("R18").Value = 3.1
$ws.Range("S18").Value = 1.83
$ws.Range("T18").Value = 1.83
$ws.Range("U18").Value = 1.2
$ws.Range("V18").Value = 4.33
$ws.Range("Y18").Value = 23
$ws.Range("Z18").Value = 34
$ws.Range("AE18").Value = 26
$ws.Range("AF18").Value = 10
$ws.Range("AI18").Value = 81
$ws.Range("AJ18").Value = 13

# Row 19
$ws.Range("G19").Value = 3.1
$ws.Range("I19").Value = 2.15
$ws.Range("K19").Value = 2.4
$ws.Range("L19").Value = 2.63
$ws.Range("S19").Value = 2.2
$ws.Range("T19").Value = 1.62
$ws.Range("W19").Value = 1.44
$ws.Range("X19").Value = 2.63
$ws.Range("AA19").Value = 12
$ws.Range("AE19").Value = 19
$ws.Range("AF19").Value = 7.5
$ws.Range("AO19").Value = 19
$ws.Range("AR19").Value = 1.83
$ws.Range("AS19").Value = 2.03

# Row 20
$ws.Range("G20").Value = 2.1
$ws.Range("H20").Value = 3.25
$ws.Range("I20").Value = 3.5
$ws.Range("J20").Value = 2.63
$ws.Range("Q20").Value = 1.73
$ws.Range("R20").Value = 2.08
$ws.Range("AD20").Value = 23
$ws.Range("AJ20").Value = 13
$ws.Range("AK20").Value = 19
$ws.Range("AL20").Value = 12

# Row 21
$ws.Range("G21").Value = 2.9
$ws.Range("H21").Value = 2.9
$ws.Range("M21").Value = 1.07
$ws.Range("N21").Value = 9
$ws.Range("O21").Value = 1.33
$ws.Range("P21").Value = 3.25
$ws.Range("Q21").Value = 2.08
$ws.Range("R21").Value = 1.73
$ws.Range("S21").Value = 3.75
$ws.Range("T21").Value = 1.25
$ws.Range("U21").Value = 1.44
$ws.Range("V21").Value = 2.63
$ws.Range("W21").Value = 1.75
$ws.Range("X21").Value = 2
$ws.Range("Y21").Value = 9.5
$ws.Range("AC21").Value = 23
$ws.Range("AD21").Value = 34
$ws.Range("AE21").Value = 9
$ws.Range("AF21").Value = 5.5
$ws.Range("AG21").Value = 13
$ws.Range("AI21").Value = 201
$ws.Range("AJ21").Value = 8.5

# Row 22
$ws.Range("N22").Value = 9
$ws.Range("O22").Value = 1.36
$ws.Range("P22").Value = 3

# Row 24
$ws.Range("G24").Value = 1.62
$ws.Range("N24").Value = 8
$ws.Range("AB24").Value = 11
$ws.Range("AE24").Value = 7
$ws.Range("AL24").Value = 21

# Row 25
$ws.Range("G25").Value = 1.78
$ws.Range("H25").Value = 3.8
$ws.Range("I25").Value = 3.9
$ws.Range("J25").Value = 2.32
$ws.Range("K25").Value = 2.3
$ws.Range("L25").Value = 4.1
$ws.Range("O25").Value = 1.19
$ws.Range("P25").Value = 4.1
$ws.Range("V25").Value = 3.15
$ws.Range("Y25").Value = 9.25
$ws.Range("Z25").Value = 10
$ws.Range("AB25").Value = 15.5
$ws.Range("AF25").Value = 7.6
$ws.Range("AJ25").Value = 15
$ws.Range("AK25").Value = 25
$ws.Range("AL25").Value = 13
$ws.Range("AM25").Value = 60
$ws.Range("AN25").Value = 30

# Row 26
$ws.Range("M26").Value = 1.06
$ws.Range("N26").Value = 10
$ws.Range("O26").Value = 1.33
$ws.Range("P26").Value = 3.25
